# Roger Tan's timesheet (week 2) - sign-off update
# - Supervisor Name is filled in with "Prakruti Sinha"
# - Supervisor signs off the sheet on row 27 ("P.S") with a sign-off date
#   of 28/02/2014 (serial 41698), matching the commit message
#   "Signed Off Time Sheets / As of 28/02/2014".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Supervisor Name value (next to "Supervisor Name:" label in row 6)
$ws.Range("G6").Value = "Prakruti Sinha"

# Supervisor signature line (row 27) + sign-off date
$ws.Range("A27").Value = "P.S"
$ws.Range("D27").Value = 41698

# Reflect the author's final selection/scroll position in the sheet view
$ws.Range("F30").Select()

Write-Output "Applied Roger Tan timesheet sign-off edits"
